$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the y-values in column C (rows 3-8) to be negative
$ws.Range("C3").Value = -1
$ws.Range("C4").Value = -1
$ws.Range("C5").Value = -2
$ws.Range("C6").Value = -2
$ws.Range("C7").Value = -2
$ws.Range("C8").Value = -2

# Update the view: scroll back to top-left and change the selection to C9
$ws.Range("A1").Select()
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C9").Select()
